$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 894.625
$ws.Range("I6").Value = 986.7143
$ws.Range("K6").Value = 2960.1429
$ws.Range("M6").Value = -2848.1429
$ws.Range("H40").Value = 4799.3335
$ws.Range("I40").Value = 2346.6667
$ws.Range("J40").Value = 8887.111000000001
$ws.Range("K40").Value = 2346.6667
$ws.Range("L40").Value = 8887.111000000001
$ws.Range("M40").Value = -2171.6667
$ws.Range("N40").Value = -9237.111000000001
$ws.Range("H132").Value = 6437
$ws.Range("I132").Value = 7546.8237
$ws.Range("K132").Value = 22640.4711
$ws.Range("M132").Value = -20110.4711
$ws.Range("H135").Value = 1794.6364
$ws.Range("J135").Value = 2753.8
$ws.Range("L135").Value = 24784.2
$ws.Range("N135").Value = -29854.2
$ws.Range("H138").Value = 5860.1
$ws.Range("I138").Value = 1963.5333
$ws.Range("J138").Value = 7158.9556
$ws.Range("K138").Value = 5890.5999
$ws.Range("L138").Value = 21476.8668
$ws.Range("M138").Value = -750.5999000000002
$ws.Range("N138").Value = -31756.8668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 142857920
$ws.Range("I4").Value = 287.8
$ws.Range("J4").Value = 500002000
$ws.Range("K4").Value = 287.8
$ws.Range("L4").Value = 500002000
$ws.Range("M4").Value = -171.8
$ws.Range("N4").Value = -500002232
$ws.Range("H5").Value = 643.7778
$ws.Range("I5").Value = 66.5
$ws.Range("K5").Value = 66.5
$ws.Range("M5").Value = 45.5
$ws.Range("H6").Value = 7599.8
$ws.Range("I6").Value = 9999
$ws.Range("K6").Value = 9999
$ws.Range("M6").Value = -9826
$ws.Range("H32").Value = 2380.423
$ws.Range("I32").Value = 2424.9607
$ws.Range("K32").Value = 2424.9607
$ws.Range("M32").Value = -2137.9607
$ws.Range("H45").Value = 19999.8
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 19999.8
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 19999.8
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -20753.8
$ws.Range("H63").Value = 7816.25
$ws.Range("I63").Value = 3500
$ws.Range("K63").Value = 3500
$ws.Range("M63").Value = -2814
$ws.Range("H66").Value = 7816.25
$ws.Range("I66").Value = 3500
$ws.Range("K66").Value = 17500
$ws.Range("M66").Value = -14068
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H122").Value = 2869.4329
$ws.Range("I122").Value = 2668.516
$ws.Range("K122").Value = 8005.548000000001
$ws.Range("M122").Value = -5555.548000000001
$ws.Range("H123").Value = 79995
$ws.Range("J123").Value = 79995
$ws.Range("L123").Value = 79995
$ws.Range("N123").Value = -89795
$ws.Range("H132").Value = 10624.523
$ws.Range("I132").Value = 4926.25
$ws.Range("J132").Value = 18222.223
$ws.Range("K132").Value = 14778.75
$ws.Range("L132").Value = 54666.66900000001
$ws.Range("M132").Value = -12248.75
$ws.Range("N132").Value = -59726.66900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 643.7778
$ws.Range("I4").Value = 66.5
$ws.Range("K4").Value = 66.5
$ws.Range("M4").Value = 48.5
$ws.Range("H19").Value = 59998
$ws.Range("J19").Value = 59998
$ws.Range("L19").Value = 59998
$ws.Range("N19").Value = -60344
$ws.Range("H26").Value = 8892.25
$ws.Range("I26").Value = 8892.25
$ws.Range("K26").Value = 8892.25
$ws.Range("M26").Value = -8600.25
$ws.Range("H35").Value = 79742.42999999999
$ws.Range("J35").Value = 79742.42999999999
$ws.Range("L35").Value = 79742.42999999999
$ws.Range("N35").Value = -80362.42999999999
$ws.Range("H82").Value = 38886.8
$ws.Range("I82").Value = 18608.5
$ws.Range("K82").Value = 18608.5
$ws.Range("M82").Value = -18225.5
$ws.Range("H85").Value = 38886.8
$ws.Range("I85").Value = 18608.5
$ws.Range("K85").Value = 18608.5
$ws.Range("M85").Value = -17282.5
$ws.Range("H134").Value = 38920.38
$ws.Range("I134").Value = 3864.65
$ws.Range("J134").Value = 116822
$ws.Range("K134").Value = 11593.95
$ws.Range("L134").Value = 350466
$ws.Range("M134").Value = -9058.950000000001
$ws.Range("N134").Value = -355536

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 586.0909
$ws.Range("I22").Value = 581.6667
$ws.Range("J22").Value = 591.4
$ws.Range("K22").Value = 581.6667
$ws.Range("L22").Value = 591.4
$ws.Range("M22").Value = -231.6667
$ws.Range("N22").Value = -1291.4
$ws.Range("H31").Value = 7245
$ws.Range("I31").Value = 3847.6667
$ws.Range("J31").Value = 8094.3335
$ws.Range("K31").Value = 3847.6667
$ws.Range("L31").Value = 8094.3335
$ws.Range("M31").Value = -3552.6667
$ws.Range("N31").Value = -8684.333500000001
$ws.Range("H34").Value = 7245
$ws.Range("I34").Value = 3847.6667
$ws.Range("J34").Value = 8094.3335
$ws.Range("K34").Value = 3847.6667
$ws.Range("L34").Value = 8094.3335
$ws.Range("M34").Value = -3645.6667
$ws.Range("N34").Value = -8498.333500000001
$ws.Range("H51").Value = 31565
$ws.Range("J51").Value = 77325
$ws.Range("L51").Value = 77325
$ws.Range("N51").Value = -78797
$ws.Range("H60").Value = 98850
$ws.Range("J60").Value = 98850
$ws.Range("L60").Value = 98850
$ws.Range("N60").Value = -99872
$ws.Range("H61").Value = 31565
$ws.Range("J61").Value = 77325
$ws.Range("L61").Value = 77325
$ws.Range("N61").Value = -78021
$ws.Range("H86").Value = 13950.091
$ws.Range("I86").Value = 6965.3335
$ws.Range("K86").Value = 6965.3335
$ws.Range("M86").Value = -5842.3335
$ws.Range("H89").Value = 13950.091
$ws.Range("I89").Value = 6965.3335
$ws.Range("K89").Value = 34826.6675
$ws.Range("M89").Value = -29210.6675
$ws.Range("H132").Value = 5991.3
$ws.Range("I132").Value = 6247.5806
$ws.Range("J132").Value = 5108.5557
$ws.Range("K132").Value = 18742.7418
$ws.Range("L132").Value = 15325.6671
$ws.Range("M132").Value = -16212.7418
$ws.Range("N132").Value = -20385.6671
$ws.Range("H134").Value = 325766.53
$ws.Range("I134").Value = 2739.842
$ws.Range("J134").Value = 837225.4399999999
$ws.Range("K134").Value = 8219.526
$ws.Range("L134").Value = 2511676.32
$ws.Range("M134").Value = -5684.526
$ws.Range("N134").Value = -2516746.32

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3201.375
$ws.Range("J131").Value = 5420.143
$ws.Range("L131").Value = 16260.429
$ws.Range("N131").Value = -26340.429
$ws.Range("H141").Value = 3715.5715
$ws.Range("I141").Value = 3834.8333
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 11504.4999
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -6324.499899999999
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5345.6665
$ws.Range("I70").Value = 5307
$ws.Range("K70").Value = 5307
$ws.Range("M70").Value = -5037
$ws.Range("H73").Value = 5345.6665
$ws.Range("I73").Value = 5307
$ws.Range("K73").Value = 5307
$ws.Range("M73").Value = -4371
$ws.Range("H124").Value = 102000
$ws.Range("J124").Value = 102000
$ws.Range("L124").Value = 102000
$ws.Range("N124").Value = -111820
$ws.Range("H128").Value = 83496.25
$ws.Range("J128").Value = 83496.25
$ws.Range("L128").Value = 83496.25
$ws.Range("N128").Value = -93456.25
$ws.Range("H130").Value = 80000
$ws.Range("J130").Value = 80000
$ws.Range("L130").Value = 80000
$ws.Range("N130").Value = -90040
$ws.Range("H132").Value = 115567.89
$ws.Range("I132").Value = 4849.5
$ws.Range("J132").Value = 204142.6
$ws.Range("K132").Value = 14548.5
$ws.Range("L132").Value = 612427.8
$ws.Range("M132").Value = -12018.5
$ws.Range("N132").Value = -617487.8
$ws.Range("H136").Value = 92657.60000000001
$ws.Range("J136").Value = 92657.60000000001
$ws.Range("L136").Value = 277972.8
$ws.Range("N136").Value = -283072.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 309357.8
$ws.Range("I122").Value = 4025.6177
$ws.Range("J122").Value = 5500005
$ws.Range("K122").Value = 12076.8531
$ws.Range("L122").Value = 16500015
$ws.Range("M122").Value = -9626.8531
$ws.Range("N122").Value = -16504915
$ws.Range("H132").Value = 4544.9556
$ws.Range("I132").Value = 3200.8215
$ws.Range("K132").Value = 9602.4645
$ws.Range("M132").Value = -7072.4645
$ws.Range("H136").Value = 5999.636
$ws.Range("I136").Value = 4599.2
$ws.Range("K136").Value = 13797.6
$ws.Range("M136").Value = -11247.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 8300
$ws.Range("I23").Value = 5001
$ws.Range("K23").Value = 5001
$ws.Range("M23").Value = -4772
$ws.Range("H81").Value = 16488.688
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 16488.688
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 19070
$ws.Range("I132").Value = 2012
$ws.Range("J132").Value = 57674.95
$ws.Range("K132").Value = 6036
$ws.Range("L132").Value = 173024.85
$ws.Range("M132").Value = -3506
$ws.Range("N132").Value = -178084.85

Write-Host "Applied all changes"